$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: target cluster changes from "sCs" to "FAPs" (self-signaling), and the
# expression metrics are recomputed accordingly.
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.876175666666667
$ws.Range("H2").Value = 5.628527
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1753453333333333
$ws.Range("N2").Value = 0.5260359999999999
$ws.Range("O2").Value = 0.414882210303281
$ws.Range("P2").Value = 0.5154067662594317
$ws.Range("Q2").Value = 0.3289786476635556
$ws.Range("R2").Value = 2.960807828972
$ws.Range("S2").Value = 0.414882210303281
$ws.Range("T2").Value = 0.5154067662594317

# Row 3: new row for the original FAPs -> sCs pairing.
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Rspo2"
$ws.Range("C3").Value = "Lgr6"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.876175666666667
$ws.Range("H3").Value = 5.628527
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2472935
$ws.Range("N3").Value = 0.494587
$ws.Range("O3").Value = 0.585117789696719
$ws.Range("P3").Value = 0.4845932337405682
$ws.Range("Q3").Value = 0.4639660472248333
$ws.Range("R3").Value = 2.783796283349
$ws.Range("S3").Value = 0.585117789696719
$ws.Range("T3").Value = 0.4845932337405682
